$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.541.28"
$ws.Range("E2").Value = "  -2.48%  "
$ws.Range("D3").Value = "2.001.83"
$ws.Range("E3").Value = "  -4.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.013"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.010"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4209"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.14"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08992"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.117"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.23%  "
$ws.Range("D13").Value = "2.013.35"
$ws.Range("E13").Value = "  -5.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.033"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.466"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001112"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06666"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.010"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.954"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.67%  "
$ws.Range("D23").Value = "29.581.43"
$ws.Range("E23").Value = "  -2.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.299"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.395"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.298"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.054"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09955"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.564"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.830"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.799"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02466"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.298"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.308"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06396"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6554"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.67"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2050"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.010"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6364"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.188"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.304"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.509"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000333"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06991"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.127"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.39%  "
